# Apply forecast/summary value updates (Optuna Attempt - go back with original)

$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("L2").Value = 0.88

# Row 3
$ws1.Range("L3").Value = 0.85

# Row 4
$ws1.Range("L4").Value = 0.97

# Row 5
$ws1.Range("D5").Value = 356
$ws1.Range("H5").Value = 10.13
$ws1.Range("L5").Value = 0.89

# Row 6
$ws1.Range("D6").Value = 319
$ws1.Range("H6").Value = 10.19
$ws1.Range("L6").Value = 0.97

# Row 7
$ws1.Range("D7").Value = 349
$ws1.Range("H7").Value = 8.4
$ws1.Range("L7").Value = 0.87

# Row 8
$ws1.Range("D8").Value = 367
$ws1.Range("H8").Value = 7.04
$ws1.Range("L8").Value = 0.9399999999999999

# Row 9
$ws1.Range("D9").Value = 334
$ws1.Range("H9").Value = 6.63
$ws1.Range("L9").Value = 1.14

# Row 10
$ws1.Range("D10").Value = 289
$ws1.Range("H10").Value = 6.51
$ws1.Range("L10").Value = 1.08

# Row 11
$ws1.Range("D11").Value = 281
$ws1.Range("H11").Value = 5.67
$ws1.Range("L11").Value = 1.18

# Row 12
$ws1.Range("D12").Value = 313
$ws1.Range("H12").Value = 4.19
$ws1.Range("L12").Value = 1.17

# Row 13
$ws1.Range("H13").Value = 2.79
$ws1.Range("L13").Value = 1.14

# Row 14
$ws1.Range("H14").Value = 1.81
$ws1.Range("L14").Value = 1.04

# Row 15
$ws1.Range("H15").Value = 0.84
$ws1.Range("L15").Value = 0.8100000000000001

# Row 16
$ws1.Range("L16").Value = 1.04

# Row 17
$ws1.Range("L17").Value = 1.09

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")

# Values in this column are stored as text, so force text format before assigning
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "5546"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "2931"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "1562"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "281"
